$wb = $excel.ActiveWorkbook

$newHandbackDate = "2016-09-03 22:46:31"
$newHandbackDateDeDe = "2016-09-03 22:46:38"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c928dcb70d3f6c39e80968f03d827e66c730d2df/e2e/2d9e27b4-2408-4729-8a10-275b1ad0eb56.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/95811c24befb4a2961d9e25441f840b4b1f790ca/e2e/2d9e27b4-2408-4729-8a10-275b1ad0eb56.md."
$latestTargetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/95811c24befb4a2961d9e25441f840b4b1f790ca/e2e/2d9e27b4-2408-4729-8a10-275b1ad0eb56.md"
$latestTargetDisplay = "2d9e27b4-2408-4729-8a10-275b1ad0eb56.md"

# ---------------- zh-cn sheet ----------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Report generated for the handback of 2d9e27b4-2408-4729-8a10-275b1ad0eb56 (row 6):
# Latest Target File (I6) becomes a hyperlink to the (newer) target markdown file.
$wsZh.Hyperlinks.Add($wsZh.Range("I6"), $latestTargetUrl, "", "", $latestTargetDisplay)

# Latest Handback File (J6): the generated xliff file name for this handback.
$wsZh.Range("J6").Value = "2d9e27b4-2408-4729-8a10-275b1ad0eb56.0270ccb1891f393f9ac4eed9e2f6d2b035cdfa0c.zh-cn.xlf"

# Latest Handback DateTime (K6): timestamp of the handback report.
$wsZh.Range("K6").Value = $newHandbackDate

# Error Detail (P6): the handback is stale relative to the latest source.
$wsZh.Range("P6").Value = $errorDetail

# Widen the Error Detail column so the message is readable.
$wsZh.Columns.Item(16).ColumnWidth = 39.17

# ---------------- de-de sheet ----------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("I6"), $latestTargetUrl, "", "", $latestTargetDisplay)

$wsDe.Range("J6").Value = "2d9e27b4-2408-4729-8a10-275b1ad0eb56.0270ccb1891f393f9ac4eed9e2f6d2b035cdfa0c.de-de.xlf"

$wsDe.Range("K6").Value = $newHandbackDateDeDe

$wsDe.Range("P6").Value = $errorDetail

$wsDe.Columns.Item(16).ColumnWidth = 39.17
